$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 7195.143
$ws.Range("I32").Value = 5993
$ws.Range("J32").Value = 7395.5
$ws.Range("K32").Value = 5993
$ws.Range("L32").Value = 7395.5
$ws.Range("M32").Value = -5667
$ws.Range("N32").Value = -8047.5
$ws.Range("H40").Value = 2761.889
$ws.Range("I40").Value = 2477.2666
$ws.Range("J40").Value = 4185
$ws.Range("K40").Value = 2477.2666
$ws.Range("L40").Value = 4185
$ws.Range("M40").Value = -2302.2666
$ws.Range("N40").Value = -4535
$ws.Range("H62").Value = 741.1667
$ws.Range("I62").Value = 689.4
$ws.Range("K62").Value = 689.4
$ws.Range("M62").Value = -65.39999999999998
$ws.Range("H64").Value = 6489.4287
$ws.Range("I64").Value = 6663.778
$ws.Range("J64").Value = 6175.6
$ws.Range("K64").Value = 6663.778
$ws.Range("L64").Value = 6175.6
$ws.Range("M64").Value = -6415.778
$ws.Range("N64").Value = -6671.6
$ws.Range("H65").Value = 741.1667
$ws.Range("I65").Value = 689.4
$ws.Range("K65").Value = 3447
$ws.Range("M65").Value = -327
$ws.Range("H67").Value = 6489.4287
$ws.Range("I67").Value = 6663.778
$ws.Range("J67").Value = 6175.6
$ws.Range("K67").Value = 6663.778
$ws.Range("L67").Value = 6175.6
$ws.Range("M67").Value = -5805.778
$ws.Range("N67").Value = -7891.6
$ws.Range("H70").Value = 5198.727
$ws.Range("I70").Value = 4898.25
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 14694.75
$ws.Range("L70").Value = 18000
$ws.Range("M70").Value = -14424.75
$ws.Range("N70").Value = -18540
$ws.Range("H73").Value = 5198.727
$ws.Range("I73").Value = 4898.25
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 14694.75
$ws.Range("L73").Value = 18000
$ws.Range("M73").Value = -13758.75
$ws.Range("N73").Value = -19872
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 45721
$ws.Range("I74").Value = 53311.086
$ws.Range("K74").Value = 53311.086
$ws.Range("M74").Value = -52437.086
$ws.Range("H77").Value = 45721
$ws.Range("I77").Value = 53311.086
$ws.Range("K77").Value = 266555.43
$ws.Range("M77").Value = -262187.43
$ws.Range("H122").Value = 2468.0952
$ws.Range("I122").Value = 2063.0908
$ws.Range("J122").Value = 2913.6
$ws.Range("K122").Value = 6189.2724
$ws.Range("L122").Value = 8740.799999999999
$ws.Range("M122").Value = -3739.2724
$ws.Range("N122").Value = -13640.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1334.3
$ws.Range("I20").Value = 1353.591
$ws.Range("K20").Value = 1353.591
$ws.Range("M20").Value = -1106.591

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2466.6667
$ws.Range("I16").Value = 2466.6667
$ws.Range("K16").Value = 2466.6667
$ws.Range("M16").Value = -2179.6667
$ws.Range("H31").Value = 1992.2439
$ws.Range("I31").Value = 1090.8667
$ws.Range("J31").Value = 4450.5454
$ws.Range("K31").Value = 1090.8667
$ws.Range("L31").Value = 4450.5454
$ws.Range("M31").Value = -795.8667
$ws.Range("N31").Value = -5040.5454
$ws.Range("H34").Value = 1992.2439
$ws.Range("I34").Value = 1090.8667
$ws.Range("J34").Value = 4450.5454
$ws.Range("K34").Value = 1090.8667
$ws.Range("L34").Value = 4450.5454
$ws.Range("M34").Value = -888.8667
$ws.Range("N34").Value = -4854.5454
$ws.Range("H62").Value = 3624.5
$ws.Range("J62").Value = 4035.3333
$ws.Range("L62").Value = 4035.3333
$ws.Range("N62").Value = -5283.3333
$ws.Range("H65").Value = 3624.5
$ws.Range("J65").Value = 4035.3333
$ws.Range("L65").Value = 20176.6665
$ws.Range("N65").Value = -26416.6665
$ws.Range("H99").Value = 4167.6875
$ws.Range("I99").Value = 4889.636
$ws.Range("J99").Value = 2579.4
$ws.Range("K99").Value = 4889.636
$ws.Range("L99").Value = 2579.4
$ws.Range("M99").Value = -3391.636
$ws.Range("N99").Value = -5575.4
$ws.Range("H107").Value = 414
$ws.Range("I107").Value = 388.18182
$ws.Range("J107").Value = 698
$ws.Range("K107").Value = 388.18182
$ws.Range("L107").Value = 698
$ws.Range("M107").Value = 1531.81818
$ws.Range("N107").Value = -4538
$ws.Range("H113").Value = 2466.6667
$ws.Range("I113").Value = 2466.6667
$ws.Range("K113").Value = 2466.6667
$ws.Range("M113").Value = -296.6667000000002
$ws.Range("H126").Value = 4167.6875
$ws.Range("I126").Value = 4889.636
$ws.Range("J126").Value = 2579.4
$ws.Range("K126").Value = 14668.908
$ws.Range("L126").Value = 7738.200000000001
$ws.Range("M126").Value = -12198.908
$ws.Range("N126").Value = -12678.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1800.5385
$ws.Range("I5").Value = 859
$ws.Range("K5").Value = 2577
$ws.Range("M5").Value = -2465
$ws.Range("H11").Value = 86
$ws.Range("I11").Value = 102.5
$ws.Range("K11").Value = 307.5
$ws.Range("M11").Value = -167.5
$ws.Range("H55").Value = 9524540
$ws.Range("I55").Value = 200
$ws.Range("K55").Value = 600
$ws.Range("M55").Value = -423
$ws.Range("H68").Value = 5666.6665
$ws.Range("I68").Value = 15000
$ws.Range("K68").Value = 45000
$ws.Range("M68").Value = -44189
$ws.Range("H71").Value = 5666.6665
$ws.Range("I71").Value = 15000
$ws.Range("K71").Value = 135000
$ws.Range("M71").Value = -130944
$ws.Range("H102").Value = 4099.2
$ws.Range("J102").Value = 5749.5
$ws.Range("L102").Value = 17248.5
$ws.Range("N102").Value = -22116.5
$ws.Range("H135").Value = 1800.5385
$ws.Range("I135").Value = 859
$ws.Range("K135").Value = 7731
$ws.Range("M135").Value = -5196

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2600.8572
$ws.Range("J97").Value = 2951.5
$ws.Range("L97").Value = 2951.5
$ws.Range("N97").Value = -3943.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5207
$ws.Range("I40").Value = 3542.5715
$ws.Range("J40").Value = 8119.75
$ws.Range("K40").Value = 3542.5715
$ws.Range("L40").Value = 8119.75
$ws.Range("M40").Value = -3406.5715
$ws.Range("N40").Value = -8391.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 13
$ws.Range("I22").Value = 13
$ws.Range("K22").Value = 13
$ws.Range("M22").Value = 280
$ws.Range("H62").Value = 9133
$ws.Range("J62").Value = 10950
$ws.Range("L62").Value = 10950
$ws.Range("N62").Value = -12198
$ws.Range("H65").Value = 9133
$ws.Range("J65").Value = 10950
$ws.Range("L65").Value = 54750
$ws.Range("N65").Value = -60990
$ws.Range("H132").Value = 5111
$ws.Range("I132").Value = 5148
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 15444
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -12914
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 2854.1738
$ws.Range("I136").Value = 1489.1082
$ws.Range("J136").Value = 8466.111000000001
$ws.Range("K136").Value = 4467.3246
$ws.Range("L136").Value = 25398.333
$ws.Range("M136").Value = -1917.3246
$ws.Range("N136").Value = -30498.333
